# Scheduled-runner style refresh of market/profit data across all craft sheets.
# For each listed row, update currentAveragePrice/NQ/HQ, LevePrice NQ/HQ and
# LeveProfit NQ/HQ columns (H:N) with freshly recalculated values. A handful of
# rows gain or lose a LeveProfit cell entirely (cleared to blank / newly populated)
# to mirror the source data exactly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 1125.25
$ws.Range("I92").Value = 892.5294
$ws.Range("J92").Value = 2444
$ws.Range("K92").Value = 892.5294
$ws.Range("L92").Value = 2444
$ws.Range("M92").Value = 355.4706
$ws.Range("N92").Value = -4940
# Row 113
$ws.Range("H113").Value = 4880.4
$ws.Range("I113").Value = 3701.25
$ws.Range("J113").Value = 5666.5
$ws.Range("K113").Value = 3701.25
$ws.Range("L113").Value = 5666.5
$ws.Range("M113").Value = -447.25
$ws.Range("N113").Value = -12174.5
# Row 141
$ws.Range("H141").Value = 8825.706
$ws.Range("I141").Value = 6986.5
$ws.Range("J141").Value = 9391.615
$ws.Range("K141").Value = 20959.5
$ws.Range("L141").Value = 28174.845
$ws.Range("M141").Value = -15779.5
$ws.Range("N141").Value = -38534.845

$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 4000
$ws.Range("I21").Value = 4000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 4000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -3626
$ws.Range("N21").Value = $null
# Row 32
$ws.Range("H32").Value = 2437.2075
$ws.Range("I32").Value = 2356.392
$ws.Range("J32").Value = 4498
$ws.Range("K32").Value = 2356.392
$ws.Range("L32").Value = 4498
$ws.Range("M32").Value = -2069.392
$ws.Range("N32").Value = -5072
# Row 61
$ws.Range("H61").Value = 2856.2104
$ws.Range("I61").Value = 2927.75
$ws.Range("J61").Value = 2474.6667
$ws.Range("K61").Value = 2927.75
$ws.Range("L61").Value = 2474.6667
$ws.Range("M61").Value = -2715.75
$ws.Range("N61").Value = -2898.6667
# Row 97
$ws.Range("H97").Value = 2575
$ws.Range("I97").Value = 2575
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2575
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2079
# Row 110
$ws.Range("H110").Value = 984.2308
$ws.Range("I110").Value = 984.2308
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 984.2308
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1060.7692
# Row 122
$ws.Range("H122").Value = 4833487
$ws.Range("I122").Value = 6175287.5
$ws.Range("J122").Value = 3005.6
$ws.Range("K122").Value = 18525862.5
$ws.Range("L122").Value = 9016.799999999999
$ws.Range("M122").Value = -18523412.5
$ws.Range("N122").Value = -13916.8
# Row 132
$ws.Range("H132").Value = 3984.3125
$ws.Range("I132").Value = 2979.3333
$ws.Range("J132").Value = 6999.25
$ws.Range("K132").Value = 8937.999899999999
$ws.Range("L132").Value = 20997.75
$ws.Range("M132").Value = -6407.999899999999
$ws.Range("N132").Value = -26057.75
# Row 136
$ws.Range("H136").Value = 2856.2104
$ws.Range("I136").Value = 2927.75
$ws.Range("J136").Value = 2474.6667
$ws.Range("K136").Value = 8783.25
$ws.Range("L136").Value = 7424.000100000001
$ws.Range("M136").Value = -6233.25
$ws.Range("N136").Value = -12524.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1771.5385
$ws.Range("I20").Value = 1890.3334
$ws.Range("J20").Value = 1504.25
$ws.Range("K20").Value = 1890.3334
$ws.Range("L20").Value = 1504.25
$ws.Range("M20").Value = -1643.3334
$ws.Range("N20").Value = -1998.25
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null
# Row 107
$ws.Range("H107").Value = 4125.943
$ws.Range("I107").Value = 1700.45
$ws.Range("J107").Value = 7359.933
$ws.Range("K107").Value = 1700.45
$ws.Range("L107").Value = 7359.933
$ws.Range("M107").Value = 219.55
$ws.Range("N107").Value = -11199.933
# Row 134
$ws.Range("H134").Value = 2247.2856
$ws.Range("I134").Value = 1553.619
$ws.Range("J134").Value = 4328.2856
$ws.Range("K134").Value = 4660.857
$ws.Range("L134").Value = 12984.8568
$ws.Range("M134").Value = -2125.857
$ws.Range("N134").Value = -18054.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
# Row 22
$ws.Range("H22").Value = 1429.7858
$ws.Range("I22").Value = 1099.8334
$ws.Range("J22").Value = 1677.25
$ws.Range("K22").Value = 1099.8334
$ws.Range("L22").Value = 1677.25
$ws.Range("M22").Value = -749.8334
$ws.Range("N22").Value = -2377.25
# Row 31
$ws.Range("H31").Value = 5410.5
$ws.Range("I31").Value = 1974.9166
$ws.Range("J31").Value = 10563.875
$ws.Range("K31").Value = 1974.9166
$ws.Range("L31").Value = 10563.875
$ws.Range("M31").Value = -1679.9166
$ws.Range("N31").Value = -11153.875
# Row 34
$ws.Range("H34").Value = 5410.5
$ws.Range("I34").Value = 1974.9166
$ws.Range("J34").Value = 10563.875
$ws.Range("K34").Value = 1974.9166
$ws.Range("L34").Value = 10563.875
$ws.Range("M34").Value = -1772.9166
$ws.Range("N34").Value = -10967.875
# Row 105
$ws.Range("H105").Value = 4508
$ws.Range("I105").Value = 6005
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 6005
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -4258
$ws.Range("N105").Value = -6505

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2571674.8
$ws.Range("I4").Value = 1422926.9
$ws.Range("J4").Value = 7166666.5
$ws.Range("K4").Value = 4268780.699999999
$ws.Range("L4").Value = 21499999.5
$ws.Range("M4").Value = -4268668.699999999
$ws.Range("N4").Value = -21500223.5
# Row 51
$ws.Range("H51").Value = 1400
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -2540
$ws.Range("N51").Value = -5420
# Row 55
$ws.Range("H55").Value = 295555.75
$ws.Range("I55").Value = 417283.4
$ws.Range("J55").Value = 3409.4
$ws.Range("K55").Value = 1251850.2
$ws.Range("L55").Value = 10228.2
$ws.Range("M55").Value = -1251673.2
$ws.Range("N55").Value = -10582.2
# Row 110
$ws.Range("H110").Value = 7341.6665
$ws.Range("I110").Value = 2027
$ws.Range("J110").Value = 9999
$ws.Range("K110").Value = 6081
$ws.Range("L110").Value = 29997
$ws.Range("M110").Value = -1991
$ws.Range("N110").Value = -38177

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 323.5
$ws.Range("I5").Value = 323.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 323.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -211.5
# Row 70
$ws.Range("H70").Value = 30312004
$ws.Range("I70").Value = 47626932
$ws.Range("J70").Value = 10874.25
$ws.Range("K70").Value = 47626932
$ws.Range("L70").Value = 10874.25
$ws.Range("M70").Value = -47626662
$ws.Range("N70").Value = -11414.25
# Row 73
$ws.Range("H73").Value = 30312004
$ws.Range("I73").Value = 47626932
$ws.Range("J73").Value = 10874.25
$ws.Range("K73").Value = 47626932
$ws.Range("L73").Value = 10874.25
$ws.Range("M73").Value = -47625996
$ws.Range("N73").Value = -12746.25
# Row 80
$ws.Range("H80").Value = 15972
$ws.Range("I80").Value = 15972
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 15972
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -14974
# Row 83
$ws.Range("H83").Value = 15972
$ws.Range("I83").Value = 15972
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 79860
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -74868
# Row 102
$ws.Range("H102").Value = 3016.8667
$ws.Range("I102").Value = 3016.8667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3016.8667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1394.8667
# Row 113
$ws.Range("H113").Value = 6151.4287
$ws.Range("I113").Value = 1576.9231
$ws.Range("J113").Value = 10116
$ws.Range("K113").Value = 1576.9231
$ws.Range("L113").Value = 10116
$ws.Range("M113").Value = 593.0769
$ws.Range("N113").Value = -14456
# Row 126
$ws.Range("H126").Value = 2944.2856
$ws.Range("I126").Value = 2916.5
$ws.Range("J126").Value = 3111
$ws.Range("K126").Value = 8749.5
$ws.Range("L126").Value = 9333
$ws.Range("M126").Value = -6279.5
$ws.Range("N126").Value = -14273
# Row 132
$ws.Range("H132").Value = 3346.739
$ws.Range("I132").Value = 3093.3809
$ws.Range("J132").Value = 6007
$ws.Range("K132").Value = 9280.1427
$ws.Range("L132").Value = 18021
$ws.Range("M132").Value = -6750.1427
$ws.Range("N132").Value = -23081

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4433.1665
$ws.Range("I7").Value = 4108.909
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 4108.909
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -3996.909
$ws.Range("N7").Value = -8224
# Row 22
$ws.Range("H22").Value = 3939.3
$ws.Range("I22").Value = 4359.8
$ws.Range("J22").Value = 3799.1333
$ws.Range("K22").Value = 4359.8
$ws.Range("L22").Value = 3799.1333
$ws.Range("M22").Value = -4064.8
$ws.Range("N22").Value = -4389.1333
# Row 27
$ws.Range("H27").Value = 3939.3
$ws.Range("I27").Value = 4359.8
$ws.Range("J27").Value = 3799.1333
$ws.Range("K27").Value = 4359.8
$ws.Range("L27").Value = 3799.1333
$ws.Range("M27").Value = -4252.8
$ws.Range("N27").Value = -4013.1333
# Row 40
$ws.Range("H40").Value = 7336.1113
$ws.Range("I40").Value = 4367.3076
$ws.Range("J40").Value = 10092.857
$ws.Range("K40").Value = 4367.3076
$ws.Range("L40").Value = 10092.857
$ws.Range("M40").Value = -4231.3076
$ws.Range("N40").Value = -10364.857
# Row 55
$ws.Range("H55").Value = 1425.5555
$ws.Range("I55").Value = 474
$ws.Range("J55").Value = 2920.8572
$ws.Range("K55").Value = 474
$ws.Range("L55").Value = 2920.8572
$ws.Range("M55").Value = -301
$ws.Range("N55").Value = -3266.8572
# Row 97
$ws.Range("H97").Value = 39537.6
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 39537.6
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 39537.6
$ws.Range("N97").Value = -41519.6
# Row 122
$ws.Range("H122").Value = 5106.7856
$ws.Range("I122").Value = 4745
$ws.Range("J122").Value = 6433.3335
$ws.Range("K122").Value = 14235
$ws.Range("L122").Value = 19300.0005
$ws.Range("M122").Value = -11785
$ws.Range("N122").Value = -24200.0005
# Row 126
$ws.Range("H126").Value = 4433.1665
$ws.Range("I126").Value = 4108.909
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 12326.727
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -9856.726999999999
$ws.Range("N126").Value = -28940
# Row 132
$ws.Range("H132").Value = 4047.9565
$ws.Range("I132").Value = 4221.1055
$ws.Range("J132").Value = 3926.111
$ws.Range("K132").Value = 12663.3165
$ws.Range("L132").Value = 11778.333
$ws.Range("M132").Value = -10133.3165
$ws.Range("N132").Value = -16838.333
# Row 136
$ws.Range("H136").Value = 4275.8096
$ws.Range("I136").Value = 4152.5293
$ws.Range("J136").Value = 4799.75
$ws.Range("K136").Value = 12457.5879
$ws.Range("L136").Value = 14399.25
$ws.Range("M136").Value = -9907.5879
$ws.Range("N136").Value = -19499.25

$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 50000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 50000
$ws.Range("N99").Value = -55990
$ws.Range("M99").Value = $null
# Row 103
$ws.Range("H103").Value = 27725
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 27725
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 27725
$ws.Range("N103").Value = -30069
# Row 107
$ws.Range("H107").Value = 9999
$ws.Range("I107").Value = 9999
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 29997
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -28077
# Row 126
$ws.Range("H126").Value = 3299.2727
$ws.Range("I126").Value = 3426.3125
$ws.Range("J126").Value = 2960.5
$ws.Range("K126").Value = 10278.9375
$ws.Range("L126").Value = 8881.5
$ws.Range("M126").Value = -7808.9375
$ws.Range("N126").Value = -13821.5
# Row 132
$ws.Range("H132").Value = 5221.7427
$ws.Range("I132").Value = 3520.5833
$ws.Range("J132").Value = 8933.362999999999
$ws.Range("K132").Value = 10561.7499
$ws.Range("L132").Value = 26800.089
$ws.Range("M132").Value = -8031.749899999999
$ws.Range("N132").Value = -31860.089
